# Apply updated Betfair Back/Lay odds to the sheet for 2025-12-04.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new numeric value
$updates = @{
    "N2"  = 5.1
    "O2"  = 1.21
    "P2"  = 2.38
    "Q2"  = 1.52
    "R2"  = 1.54
    "S2"  = 2.26
    "T2"  = 1.54
    "U2"  = 2.5
    "X2"  = 28
    "Y2"  = 15.5
    "Z2"  = 20
    "AB2" = 20
    "AC2" = 10
    "AD2" = 13
    "AE2" = 29
    "AF2" = 28
    "AG2" = 14
    "AH2" = 16
    "AK2" = 36
    "AN2" = 19.5
    "AO2" = 14.5

    "H3"  = 1.4

    "P4"  = 1.63
    "Q4"  = 2.1

    "F5"  = 7.4
    "G5"  = 10.5
    "H5"  = 1.49
    "I5"  = 1.6
    "J5"  = 3.95
    "P5"  = 1.81
    "Q5"  = 1.98

    "F6"  = 4.2
    "G6"  = 8.199999999999999
    "H6"  = 1.55
    "I6"  = 1.78
    "J6"  = 3.65
    "K6"  = 5.4
    "P6"  = 2.1
    "Q6"  = 1.73

    "S8"  = 2.42
    "AO8" = 85

    "F10"  = 2.06
    "I10"  = 4.3
    "K10"  = 3.7
    "L10"  = 1.45
    "Q10"  = 2.08
    "V10"  = 1.3
    "X10"  = 15
    "AC10" = 8
    "AF10" = 16
    "AH10" = 24
    "AM10" = 150
    "AN10" = 24

    "F11" = 1.84
    "G11" = 2.16
    "H11" = 3.9
    "I11" = 6.6
    "J11" = 3.15
    "K11" = 5.1
    "P11" = 1.24
    "Q11" = 1.01
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
